$p = $ppt.ActivePresentation

# --- Slide 2: "Users can see this" / "Or this" textboxes -----------------
$s2 = $p.Slides.Item(2)

# "Users can see this" -> "Users can see this – easy to read"
# Shape moves up/left and grows taller (autofit handles the height).
$tb5 = $s2.Shapes.Item("TextBox 5")
$tb5.Left = 23.617087364196983
$tb5.Top = 64.01732254028369
$tb5.TextFrame.TextRange.Text = "Users can see this " + [char]0x2013 + " easy to read"

# "Or this" -> "Or this – underlying database tables and columns "
# Shape shifts 1 EMU left and grows much wider; set size before text so
# autofit computes the (unchanged) height against the final width.
$tb6 = $s2.Shapes.Item("TextBox 6")
$tb6.Left = 288.4444122314455
$tb6.Width = 446.6939239501954
$tb6.TextFrame.TextRange.Text = "Or this " + [char]0x2013 + " underlying database tables and columns "

# --- Slide 4: Import mode bullet list -------------------------------------
$s4 = $p.Slides.Item(4)
$tb6_s4 = $s4.Shapes.Item("TextBox 6")
$tf = $tb6_s4.TextFrame
$tr = $tf.TextRange
$limitPara = $tr.Paragraphs(4)
$limitPara.InsertAfter("`rLarger with Premium")
